$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-5
# from serial 45185 (2023-09-16) to serial 45204 (2023-10-05)
$ws.Range("C2").Value2 = 45204
$ws.Range("C3").Value2 = 45204
$ws.Range("C4").Value2 = 45204
$ws.Range("C5").Value2 = 45204
